$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.921.64"
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("D3").Value = "'1.637.37"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "'214.63"
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("E8").Value = "  -0.91%  "
$ws.Range("E9").Value = "  +0.55%  "
$ws.Range("D10").Value = "'19.60"
$ws.Range("E10").Value = "  -0.37%  "
$ws.Range("D11").Value = "'0.0793"
$ws.Range("E11").Value = "  +0.29%  "
$ws.Range("D12").Value = "'1.863.02"
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "'1.658.94"
$ws.Range("E13").Value = "  +1.05%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'4.24"
$ws.Range("E14").Value = "  -0.39%  "
$ws.Range("E15").Value = "  -1.53%  "
$ws.Range("E16").Value = "  -0.65%  "
$ws.Range("D17").Value = "'62.57"
$ws.Range("E17").Value = "  -0.61%  "
$ws.Range("D18").Value = "'25.939.68"
$ws.Range("E18").Value = "  +0.14%  "
$ws.Range("E19").Value = "  +0.20%  "
$ws.Range("D20").Value = "'193.79"
$ws.Range("E20").Value = "  +1.06%  "
$ws.Range("E21").Value = "  -1.22%  "
$ws.Range("E22").Value = "  -0.64%  "
$ws.Range("E23").Value = "  -1.00%  "
$ws.Range("E24").Value = "  +0.52%  "
$ws.Range("D25").Value = "'143.82"
$ws.Range("E25").Value = "  +1.09%  "
$ws.Range("E26").Value = "  +0.14%  "
$ws.Range("D27").Value = "'0.127"
$ws.Range("E27").Value = "  +3.05%  "
$ws.Range("D28").Value = "'6.84"
$ws.Range("E28").Value = "  -0.33%  "
$ws.Range("E29").Value = "  -0.70%  "
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("D31").Value = "'0.0501"
$ws.Range("E31").Value = "  +1.58%  "
$ws.Range("D32").Value = "'3.30"
$ws.Range("E32").Value = "  -1.09%  "
$ws.Range("D33").Value = "'3.22"
$ws.Range("E33").Value = "  -0.98%  "
$ws.Range("E34").Value = "  -2.46%  "
$ws.Range("E35").Value = "  +1.35%  "
$ws.Range("D36").Value = "'0.903"
$ws.Range("E36").Value = "  -0.56%  "
$ws.Range("D37").Value = "'1.137.87"
$ws.Range("E37").Value = "  -0.83%  "
$ws.Range("D38").Value = "'0.545"
$ws.Range("E38").Value = "  -0.08%  "
$ws.Range("D39").Value = "'2.47"
$ws.Range("E39").Value = "  -1.48%  "
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("E41").Value = "  +0.13%  "
$ws.Range("D42").Value = "'99.41"
$ws.Range("E42").Value = "  -1.40%  "
$ws.Range("D43").Value = "'0.799"
$ws.Range("E43").Value = "  -0.53%  "
$ws.Range("E44").Value = "  -3.89%  "
$ws.Range("D45").Value = "'1.772.45"
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("E46").Value = "  +3.85%  "
$ws.Range("D47").Value = "'56.43"
$ws.Range("E47").Value = "  +1.47%  "
$ws.Range("D48").Value = "'0.0530"
$ws.Range("E48").Value = "  +3.56%  "
$ws.Range("D49").Value = "'1.45"
$ws.Range("E49").Value = "  -0.85%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'7.66"
$ws.Range("E50").Value = "  +0.82%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").Value = "'0.415"
$ws.Range("E51").Value = "  -0.43%  "
